$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "ALU op" mini-table (was rows 2-12 cols M-Q) down to rows 13-23 cols D-H ---

# Header "ALU op" merged banner: M2:Q2 -> D13:H13
$ws.Range("M2:Q2").UnMerge()
$ws.Range("M2:Q2").Cut($ws.Range("D13:H13"))
$ws.Range("M2:Q2").Clear()

# Bit index header row: L3:P3 -> C14:G14
$ws.Range("L3:P3").Cut($ws.Range("C14:G14"))
$ws.Range("M3:P3").Clear()

# Data rows: M{row}:Q{row} -> D{row+11}:H{row+11}
$ws.Range("M4:Q4").Cut($ws.Range("D15:H15"))
$ws.Range("M4:Q4").Clear()
$ws.Range("M5:Q5").Cut($ws.Range("D16:H16"))
$ws.Range("M5:Q5").Clear()
$ws.Range("M6:Q6").Cut($ws.Range("D17:H17"))
$ws.Range("M6:Q6").Clear()
$ws.Range("M7:Q7").Cut($ws.Range("D18:H18"))
$ws.Range("M7:Q7").Clear()
$ws.Range("M8:Q8").Cut($ws.Range("D19:H19"))
$ws.Range("M8:Q8").Clear()
$ws.Range("M9:Q9").Cut($ws.Range("D20:H20"))
$ws.Range("M9:Q9").Clear()
$ws.Range("M10:Q10").Cut($ws.Range("D21:H21"))
$ws.Range("M10:Q10").Clear()
$ws.Range("M11:Q11").Cut($ws.Range("D22:H22"))
$ws.Range("M11:Q11").Clear()
$ws.Range("M12:Q12").Cut($ws.Range("D23:H23"))
$ws.Range("M12:Q12").Clear()

# --- Move the "Ld/str" mini-table (was rows 15-18 cols N-Q) down to rows 26-29 cols E-H ---
$ws.Range("N15:P15").Cut($ws.Range("E26:G26"))
$ws.Range("N15:P15").Clear()
$ws.Range("O16:P16").Cut($ws.Range("F27:G27"))
$ws.Range("O16:P16").Clear()
$ws.Range("P17:Q17").Cut($ws.Range("G28:H28"))
$ws.Range("P17:Q17").Clear()
$ws.Range("P18:Q18").Cut($ws.Range("G29:H29"))
$ws.Range("P18:Q18").Clear()

# --- Move the "Tru/fal" mini-table (was rows 20-23 cols N-Q) down to rows 31-34 cols E-H ---
$ws.Range("N20:P20").Cut($ws.Range("E31:G31"))
$ws.Range("N20:P20").Clear()
$ws.Range("O21:P21").Cut($ws.Range("F32:G32"))
$ws.Range("O21:P21").Clear()
$ws.Range("P22:Q22").Cut($ws.Range("G33:H33"))
$ws.Range("P22:Q22").Clear()
$ws.Range("P23:Q23").Cut($ws.Range("G34:H34"))
$ws.Range("P23:Q23").Clear()

# --- Fix up merges: add the new D13:H13 merge ---
$ws.Range("D13:H13").Merge()

# --- New "Notes" column (L) ---
$ws.Range("L3").Value = "Notes"
$ws.Range("L4").Value = "Pushes to top of stack"
$ws.Range("L5").Value = "Calculates top +-*<<>> next"
$ws.Range("L7").Value = "Loads/stores from memory location in top"
$ws.Range("L6").Value = "Pops from top of stack N amount"
$ws.Range("L8").Value = "Jumps to location in top if next is true/false"

$ws.Range("L7").HorizontalAlignment = -4131
$ws.Range("L7").VerticalAlignment = -4108
$ws.Range("L7").Borders(7).LineStyle = 1

# Column L width (Notes column)
$ws.Columns.Item(12).ColumnWidth = 38.6

# Selection left where the author's cursor ended up
$ws.Range("L9").Select()

Write-Output "done"
